$d = $word.ActiveDocument

# The title paragraph currently reads "Documentation : Website B" as a
# single run, followed by the (hidden) _GoBack bookmark. We need to turn
# it into three runs: "Documentation :" / " Level 1" / " Website B",
# with the bookmark sitting between the 2nd and 3rd runs (i.e. right
# after " Level 1", same place it sits today relative to the end of the
# original text).

$p1 = $d.Paragraphs.Item(1)
$fullText = $p1.Range.Text
$marker = " Website B"
$idx = $fullText.IndexOf($marker)

$paraStart = $p1.Range.Start
$insertAt = $paraStart + $idx
$insertEnd = $insertAt + $marker.Length

# Insert " Level 1" right before " Website B" (and therefore right
# before the bookmark, which currently sits at the end of the run).
$tail = $d.Range($insertAt, $insertEnd)
$tail.InsertBefore(" Level 1")

# Force Word to split the merged run into separate runs by nudging a
# character formatting property on just the " Level 1" span, then
# restoring it back to the paragraph's actual value (bold). This keeps
# " Level 1" / " Website B" as independent <w:r> runs instead of being
# re-coalesced into a single run.
$newLevelStart = $insertAt
$newLevelEnd = $insertAt + (" Level 1").Length
$levelRange = $d.Range($newLevelStart, $newLevelEnd)
$levelRange.Bold = 0

$levelRange2 = $d.Range($newLevelStart, $newLevelEnd)
$levelRange2.Bold = 1
